$d = $word.ActiveDocument

$pairs = @(
    @("94÷4=", "74÷9="),
    @("18÷9=", "99÷2="),
    @("47÷4=", "31÷4="),
    @("49÷7=", "80÷9="),
    @("74÷7=", "32÷5="),
    @("36÷3=", "89÷7="),
    @("66÷6=", "76÷9="),
    @("27÷2=", "84÷6="),
    @("56÷3=", "23÷8="),
    @("89÷4=", "14÷9="),
    @("91÷8=", "11÷2="),
    @("91÷4=", "45÷6="),
    @("75÷5=", "18÷2="),
    @("91÷2=", "16÷6="),
    @("53÷7=", "61÷3="),
    @("38÷2=", "31÷8="),
    @("13÷3=", "51÷4="),
    @("27÷8=", "67÷6="),
    @("38÷6=", "59÷9="),
    @("10÷4=", "36÷6="),
    @("45÷8=", "99÷8="),
    @("57÷3=", "38÷7="),
    @("44÷2=", "62÷8="),
    @("51÷6=", "93÷7="),
    @("44÷9=", "41÷9=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
